$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell C10 (the "From" value for rule R30) from 18 to 100
$ws.Range("C10").Value = 100

$wb.Save()
